# Apply "Only product type 1-2 and 2-3" change.
# Sheet "Ark3" (3rd sheet) currently has a row for ChrisKunden with product
# type "fisk3til4" (3-4). That product type should be removed and replaced
# by rows for the two allowed product types (fisk1til2 and fisk2til3), in
# line with the other customer (KristineKunden) who already has both.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Ark3")

# Row 4 used to reference "fisk3til4" - change it to "fisk1til2"
$ws3.Range("B4").Value = "fisk1til2"

# Add a new row 5 for ChrisKunden with "fisk2til3"
$ws3.Range("A5").Value = "ChrisKunden"
$ws3.Range("B5").Value = "fisk2til3"
$ws3.Range("C5").Value = 3

# Update selections / active cells to reflect the editing session
$ws1 = $wb.Worksheets.Item("Ark1")
$ws1.Range("D7").Select()

$ws3.Select()
$ws3.Range("C10").Select()

$wb.Activate()
